$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 5 de Mayo de 2020 a las 12:03"

# Row 16: Belgica
$ws.Cells.Item(16, 1).Value = "Belgica"
$ws.Cells.Item(16, 2).Value = 50509
$ws.Cells.Item(16, 3).Value = 242
$ws.Cells.Item(16, 4).Value = 12441
$ws.Cells.Item(16, 5).Value = 30052
$ws.Cells.Item(16, 6).Value = 646
$ws.Cells.Item(16, 7).Value = 92
$ws.Cells.Item(16, 8).Value = 8016

# Row 21: Suiza
$ws.Cells.Item(21, 1).Value = "Suiza"
$ws.Cells.Item(21, 2).Value = 30009
$ws.Cells.Item(21, 3).Value = 28
$ws.Cells.Item(21, 4).Value = 25200
$ws.Cells.Item(21, 5).Value = 3025
$ws.Cells.Item(21, 6).Value = 141
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 1784

# Row 39: Indonesia
$ws.Cells.Item(39, 1).Value = "Indonesia"
$ws.Cells.Item(39, 2).Value = 12071
$ws.Cells.Item(39, 3).Value = 484
$ws.Cells.Item(39, 4).Value = 2197
$ws.Cells.Item(39, 5).Value = 9002
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 8
$ws.Cells.Item(39, 8).Value = 872

# Row 42: Dinamarca
$ws.Cells.Item(42, 1).Value = "Dinamarca"
$ws.Cells.Item(42, 2).Value = 9821
$ws.Cells.Item(42, 3).Value = 151
$ws.Cells.Item(42, 4).Value = 7088
$ws.Cells.Item(42, 5).Value = 2240
$ws.Cells.Item(42, 6).Value = 57
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 493

# Row 43: Filipinas
$ws.Cells.Item(43, 1).Value = "Filipinas"
$ws.Cells.Item(43, 2).Value = 9684
$ws.Cells.Item(43, 3).Value = 199
$ws.Cells.Item(43, 4).Value = 1408
$ws.Cells.Item(43, 5).Value = 7639
$ws.Cells.Item(43, 6).Value = 31
$ws.Cells.Item(43, 7).Value = 14
$ws.Cells.Item(43, 8).Value = 637

# Row 54: Finlandia
$ws.Cells.Item(54, 1).Value = "Finlandia"
$ws.Cells.Item(54, 2).Value = 5412
$ws.Cells.Item(54, 3).Value = 85
$ws.Cells.Item(54, 4).Value = 3500
$ws.Cells.Item(54, 5).Value = 1672
$ws.Cells.Item(54, 6).Value = 49
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 240

# Row 59: Moldavia
$ws.Cells.Item(59, 1).Value = "Moldavia"
$ws.Cells.Item(59, 2).Value = 4248
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 1544
$ws.Cells.Item(59, 5).Value = 2571
$ws.Cells.Item(59, 6).Value = 237
$ws.Cells.Item(59, 7).Value = 1
$ws.Cells.Item(59, 8).Value = 133

# Row 60: Kazajistan
$ws.Cells.Item(60, 1).Value = "Kazajistan"
$ws.Cells.Item(60, 2).Value = 4160
$ws.Cells.Item(60, 3).Value = 111
$ws.Cells.Item(60, 4).Value = 1207
$ws.Cells.Item(60, 5).Value = 2924
$ws.Cells.Item(60, 6).Value = 40
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 29

# Row 63: Afganistan
$ws.Cells.Item(63, 1).Value = "Afganistan"
$ws.Cells.Item(63, 2).Value = 3224
$ws.Cells.Item(63, 3).Value = 330
$ws.Cells.Item(63, 4).Value = 421
$ws.Cells.Item(63, 5).Value = 2708
$ws.Cells.Item(63, 6).Value = 7
$ws.Cells.Item(63, 7).Value = 5
$ws.Cells.Item(63, 8).Value = 95

# Row 64: Hungria
$ws.Cells.Item(64, 1).Value = "Hungria"
$ws.Cells.Item(64, 2).Value = 3065
$ws.Cells.Item(64, 3).Value = 30
$ws.Cells.Item(64, 4).Value = 709
$ws.Cells.Item(64, 5).Value = 1993
$ws.Cells.Item(64, 6).Value = 55
$ws.Cells.Item(64, 7).Value = 12
$ws.Cells.Item(64, 8).Value = 363

# Row 65: Tailandia
$ws.Cells.Item(65, 1).Value = "Tailandia"
$ws.Cells.Item(65, 2).Value = 2988
$ws.Cells.Item(65, 3).Value = 1
$ws.Cells.Item(65, 4).Value = 2747
$ws.Cells.Item(65, 5).Value = 187
$ws.Cells.Item(65, 6).Value = 61
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 54

# Row 72: Uzbekistan
$ws.Cells.Item(72, 1).Value = "Uzbekistan"
$ws.Cells.Item(72, 2).Value = 2189
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 1454
$ws.Cells.Item(72, 5).Value = 725
$ws.Cells.Item(72, 6).Value = 8
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 10

# Row 92: Hong Kong
$ws.Cells.Item(92, 1).Value = "Hong Kong"
$ws.Cells.Item(92, 2).Value = 1041
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 920
$ws.Cells.Item(92, 5).Value = 117
$ws.Cells.Item(92, 6).Value = 1
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 4

# Row 130: Montenegro
$ws.Cells.Item(130, 1).Value = "Montenegro"
$ws.Cells.Item(130, 2).Value = 324
$ws.Cells.Item(130, 3).Value = 1
$ws.Cells.Item(130, 4).Value = 253
$ws.Cells.Item(130, 5).Value = 63
$ws.Cells.Item(130, 6).Value = 2
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 8

# Row 144: Etiopia
$ws.Cells.Item(144, 1).Value = "Etiopia"
$ws.Cells.Item(144, 2).Value = 145
$ws.Cells.Item(144, 3).Value = 5
$ws.Cells.Item(144, 4).Value = 91
$ws.Cells.Item(144, 5).Value = 50
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 4

# Row 145: Gibraltar
$ws.Cells.Item(145, 1).Value = "Gibraltar"
$ws.Cells.Item(145, 2).Value = 144
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 133
$ws.Cells.Item(145, 5).Value = 11
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 0
